# "update image and favicon"
#
# 1. Bump the cached "datetimeFigureOut" date field text from 7/17/24 to
#    7/18/24 everywhere it appears (slide master + every slide layout).
# 2. Re-style/re-layout the favicon artwork on slide 1:
#      - the big circle (Oval 4) moves slightly and switches from the
#        accent6 theme color to a flat srgbClr 00B0F0 fill
#      - the tallest bar (Rectangle 5) shifts position
#      - the "Rectangle 12" bar is removed entirely
#      - the remaining two bars (old "Rectangle 13" / "Rectangle 14") are
#        resized/repositioned to fill the gap

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Date placeholder text: 7/17/24 -> 7/18/24, on the slide master and
#    on every slide layout.
# ---------------------------------------------------------------------
function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
                if ($sh.TextFrame.TextRange.Text -eq "7/17/24") {
                    $sh.TextFrame.TextRange.Text = "7/18/24"
                }
            }
        }
    }
}

Update-DatePlaceholder $p.SlideMaster.Shapes
for ($li = 1; $li -le $p.SlideMaster.CustomLayouts.Count; $li++) {
    Update-DatePlaceholder $p.SlideMaster.CustomLayouts.Item($li).Shapes
}

# ---------------------------------------------------------------------
# 2) Slide 1 favicon artwork.
# ---------------------------------------------------------------------
$s = $p.Slides.Item(1)

# -- Oval 4: nudge position, flat-color fill instead of theme accent6 --
$oval = $s.Shapes.Item("Oval 4")
$oval.Left = 103845 / 12700
$oval.Top = 103845 / 12700
$oval.Fill.ForeColor.RGB = 15773696   # COM RGB() packs as B*65536+G*256+R -> 00B0F0

# -- Rectangle 5: reposition only (size unchanged) --
$rect5 = $s.Shapes.Item("Rectangle 5")
$rect5.Left = 2082798 / 12700
$rect5.Top = 31298 / 12700

# -- Remove the "Rectangle 12" bar entirely --
$s.Shapes.Item("Rectangle 12").Delete()

# -- Rectangle 13: reposition + resize to close the gap --
$rect13 = $s.Shapes.Item("Rectangle 13")
$rect13.Left = 1667646 / 12700
$rect13.Top = 1318190 / 12700
$rect13.Width = 514353 / 12700
$rect13.Height = 2271793 / 12700

# -- Rectangle 14: reposition only (size unchanged) --
$rect14 = $s.Shapes.Item("Rectangle 14")
$rect14.Left = 1182326 / 12700
$rect14.Top = 2675250 / 12700
